$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Match header style of existing headers (copy format from A1, which already
# carries the bold/centered/bordered header style) without disturbing values
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Boolean outlier flag values for rows 2-12
$values = @(
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,0),
    @(0,0,0),
    @(0,0,0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    $ws.Range("F$row").Value = [bool]$rowVals[0]
    $ws.Range("G$row").Value = [bool]$rowVals[1]
    $ws.Range("H$row").Value = [bool]$rowVals[2]
}
